$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 4-11: replace the "2012-0" text value (style s="9", shared string) with
# the plain number 2012, with no explicit cell style (defaults to General).
foreach ($r in 4..11) {
    $cell = $ws.Range("D$r")
    $cell.Style = "Normal"
    $cell.Value = 2012
}

# Rows 15-22: these cells were blank (style s="6") and become the text
# "2012-0", using the same style (numFmtId "0.00") that D4:D11 used to have.
foreach ($r in 15..22) {
    $cell = $ws.Range("D$r")
    $cell.NumberFormat = "0.00"
    $cell.Value = "2012-0"
}
